$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (Subj ids) - only B1:E1 changed
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) - B2:E2 changed
$ws.Range("B2").Value = 590.31889494079167
$ws.Range("C2").Value = 401.8965782137555
$ws.Range("D2").Value = 755.62365799511065
$ws.Range("E2").Value = 383.60313119269222

# Row 3 (STR) - B3:E3 changed
$ws.Range("B3").Value = 639.46245178292122
$ws.Range("C3").Value = 404.9306028098606
$ws.Range("D3").Value = 908.20732508596382
$ws.Range("E3").Value = 490.67665434144442

# Update the selection to match the new sqref B1:E3
$ws.Range("B1:E3").Select()
